$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 7-11 (A7..A11 correspond to Q6..Q10 labels already present)
# with new forecast error statistics values.

$ws.Range("B7").Value = 0.1732741012286916
$ws.Range("C7").Value = 1.574991177320644
$ws.Range("D7").Value = 9.359491832855587
$ws.Range("E7").Value = 3.05932865721478
$ws.Range("F7").Value = 3.095418501638419
$ws.Range("G7").Value = 38

$ws.Range("B8").Value = 0.2044873955334486
$ws.Range("C8").Value = 1.562635388531142
$ws.Range("D8").Value = 9.505286999794503
$ws.Range("E8").Value = 3.083064546809636
$ws.Range("F8").Value = 3.118709047737856
$ws.Range("G8").Value = 37

$ws.Range("B9").Value = -0.09975482970692948
$ws.Range("C9").Value = 2.241897685767165
$ws.Range("D9").Value = 15.38488210364042
$ws.Range("E9").Value = 3.922356702754152
$ws.Range("F9").Value = 4.022951401199263
$ws.Range("G9").Value = 20

$ws.Range("B10").Value = -0.8689250956944341
$ws.Range("C10").Value = 2.089385287080378
$ws.Range("D10").Value = 12.84994429113759
$ws.Range("E10").Value = 3.584681895390105
$ws.Range("F10").Value = 3.619782165404088
$ws.Range("G10").Value = 13

$ws.Range("B11").Value = 0.9642800447058392
$ws.Range("C11").Value = 1.347199366044221
$ws.Range("D11").Value = 3.741412787431818
$ws.Range("E11").Value = 1.934273193587663
$ws.Range("F11").Value = 1.874692235679607
$ws.Range("G11").Value = 5
